$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current values from A2:A46 into an array
$values = @()
for ($r = 2; $r -le 46; $r++) {
    $values += $ws.Cells.Item($r, 1).Value2
}

# Remove the two "bad zombie" values (74 and 142) from their current
# positions and append them at the end of the list.
$remaining = @()
foreach ($v in $values) {
    if ($v -ne 74 -and $v -ne 142) {
        $remaining += $v
    }
}
$remaining += 74
$remaining += 142

# Write the reordered values back to A2:A46
for ($i = 0; $i -lt $remaining.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $remaining[$i]
}
